$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New course headers (L1, U1) mirroring existing C1 ("Engenharia de Computação") ---
# Merge first, then set value/format on the (now single) anchor cell - doing it in
# this order keeps the format change scoped to the anchor cell only.
$ws.Range("L1:N1").Merge()
$ws.Range("U1:W1").Merge()

$ws.Range("L1").Value = "Engenharia Elétrica"
$ws.Range("U1").Value = "Engenharia Mecânica"

$ws.Range("L1").HorizontalAlignment = -4108
$ws.Range("L1").VerticalAlignment = -4108
$ws.Range("U1").HorizontalAlignment = -4108
$ws.Range("U1").VerticalAlignment = -4108

# --- Duplicate the 10 schedule tables (periods 1-10) into columns J:P and S:Y ---
$blockStarts = @(2,18,34,50,66,82,98,114,130,146)

foreach ($r in $blockStarts) {
    # Header row: columns A-G (1-7) -> J-P (10-16) and S-Y (19-25)
    for ($c = 1; $c -le 7; $c++) {
        $src = $ws.Cells.Item($r, $c)
        $dst1 = $ws.Cells.Item($r, $c + 9)
        $dst2 = $ws.Cells.Item($r, $c + 18)
        $src.Copy($dst1)
        $src.Copy($dst2)
    }
    # Data rows: column A (1) -> column J (10) and column S (19)
    for ($rr = $r + 1; $rr -le $r + 14; $rr++) {
        $src = $ws.Cells.Item($rr, 1)
        $dst1 = $ws.Cells.Item($rr, 10)
        $dst2 = $ws.Cells.Item($rr, 19)
        $src.Copy($dst1)
        $src.Copy($dst2)
    }
}

# --- Match column widths of column A (bestFit) on the two new leading columns ---
$ws.Columns.Item(10).ColumnWidth = 10.8
$ws.Columns.Item(19).ColumnWidth = 10.8

Write-Host "done"
